$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Actualiza base de datos EC: nuevos periodos de mora (parte 1 de nuevos estados de cuenta)
$ws.Range("E16").Value = "2403"
$ws.Range("E17").Value = "2404"
$ws.Range("E18").Value = "2405"

# Actualiza los valores de mora
$ws.Range("G16").Value = 1923000
$ws.Range("G17").Value = 1923000
$ws.Range("G18").Value = 1923000
